# The checklist item "- date " is removed: the feature it tracked now
# works (per the commit message "la date fonctionne !"), so the TODO
# line is cleared out, leaving the paragraph blank (same as the
# surrounding blank-line paragraphs).

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("- date ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
